# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   on every sheet (Overview zh-cn/de-de "Priority"+"Content Duplicate" cols,
#   and the per-language "Status" column).
# - The per-language sheets now carry a real "Latest Target File" (hyperlinked
#   to the source .md) and "Latest Handback File" (the generated .xlf), plus
#   an updated "Latest Handback DateTime".
# - A few columns are widened to fit the new, longer values.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdFile   = "cdf759f2-9468-423e-a693-951d82e646ff.md"
$mdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edfdd6c9c72f7b8fb462c8a4b118a5b1cfba6838/e2e/cdf759f2-9468-423e-a693-951d82e646ff.md"

# Excel's ColumnWidth (character units) gets re-quantized to pixel-snapped
# widths on save, so feed it values that land on the closest achievable
# raw width to the target layout (~30 and exactly 40 character-widths).
$wideColWidth   = 29.17
$widerColWidth  = 39.14

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusNew
$overview.Range("F2").Value = $statusNew
$overview.Columns.Item(5).ColumnWidth = $wideColWidth
$overview.Columns.Item(6).ColumnWidth = $wideColWidth

# ---- Per-language sheets (zh-cn / de-de) ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusNew
$zh.Columns.Item(3).ColumnWidth = $wideColWidth
$zh.Columns.Item(9).ColumnWidth = $widerColWidth
$zh.Columns.Item(10).ColumnWidth = $widerColWidth
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl, "", "", $mdFile)
$zh.Range("J2").Value = "cdf759f2-9468-423e-a693-951d82e646ff.d562e386f61f40ad0d449cc87e94029e48370a84.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-23 02:57:27"

$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusNew
$de.Columns.Item(3).ColumnWidth = $wideColWidth
$de.Columns.Item(9).ColumnWidth = $widerColWidth
$de.Columns.Item(10).ColumnWidth = $widerColWidth
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl, "", "", $mdFile)
$de.Range("J2").Value = "cdf759f2-9468-423e-a693-951d82e646ff.d562e386f61f40ad0d449cc87e94029e48370a84.de-de.xlf"
$de.Range("K2").Value = "2016-08-23 02:57:34"

Write-Output "Handback report generated"
